# Update scripts with new TPM-derived NATMI statistics.
# Ligand-expressing-cells (E:J) and Receptor-expressing-cells (K:T) columns
# were recomputed after the TPM update; write the refreshed values in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04983966666666667
$ws.Range("H2").Value = 0.149519
$ws.Range("I2").Value = 0.1823731600337622
$ws.Range("J2").Value = 0.1823731600337622
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2023976666666667
$ws.Range("N2").Value = 0.6071930000000001
$ws.Range("O2").Value = 0.03663970451354832
$ws.Range("P2").Value = 0.03663970451354832
$ws.Range("Q2").Value = 0.01008743224077778
$ws.Range("R2").Value = 0.09078689016700002
$ws.Range("S2").Value = 0.006682098694839107
$ws.Range("T2").Value = 0.006682098694839106
$ws.Range("O3").Value = 0.08641717548188978
$ws.Range("P3").Value = 0.08641717548188979
$ws.Range("Q3").Value = 0.02379187861055556
$ws.Range("R3").Value = 0.214126907495
$ws.Range("S3").Value = 0.0157601733738244
$ws.Range("T3").Value = 0.0157601733738244
$ws.Range("O4").Value = 0.876943120004562
$ws.Range("P4").Value = 0.876943120004562
$ws.Range("Q4").Value = 0.2414349247492222
$ws.Range("R4").Value = 2.172914322743
$ws.Range("S4").Value = 0.1599308879650987
$ws.Range("T4").Value = 0.1599308879650987
$ws.Range("I5").Value = 0.4031422744592926
$ws.Range("J5").Value = 0.4031422744592926
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2023976666666667
$ws.Range("N5").Value = 0.6071930000000001
$ws.Range("O5").Value = 0.03663970451354832
$ws.Range("P5").Value = 0.03663970451354832
$ws.Range("Q5").Value = 0.02229862319788889
$ws.Range("R5").Value = 0.200687608781
$ws.Range("S5").Value = 0.01477101381310828
$ws.Range("T5").Value = 0.01477101381310828
$ws.Range("I6").Value = 0.4031422744592926
$ws.Range("J6").Value = 0.4031422744592926
$ws.Range("O6").Value = 0.08641717548188978
$ws.Range("P6").Value = 0.08641717548188979
$ws.Range("S6").Value = 0.03483841667611686
$ws.Range("T6").Value = 0.03483841667611686
$ws.Range("I7").Value = 0.4031422744592926
$ws.Range("J7").Value = 0.4031422744592926
$ws.Range("O7").Value = 0.876943120004562
$ws.Range("P7").Value = 0.876943120004562
$ws.Range("S7").Value = 0.3535328439700675
$ws.Range("T7").Value = 0.3535328439700675
$ws.Range("I8").Value = 0.4144845655069452
$ws.Range("J8").Value = 0.4144845655069451
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2023976666666667
$ws.Range("N8").Value = 0.6071930000000001
$ws.Range("O8").Value = 0.03663970451354832
$ws.Range("P8").Value = 0.03663970451354832
$ws.Range("Q8").Value = 0.02292598849866667
$ws.Range("R8").Value = 0.206333896488
$ws.Range("S8").Value = 0.01518659200560093
$ws.Range("T8").Value = 0.01518659200560093
$ws.Range("I9").Value = 0.4144845655069452
$ws.Range("J9").Value = 0.4144845655069451
$ws.Range("O9").Value = 0.08641717548188978
$ws.Range("P9").Value = 0.08641717548188979
$ws.Range("S9").Value = 0.03581858543194852
$ws.Range("T9").Value = 0.03581858543194852
$ws.Range("I10").Value = 0.4144845655069452
$ws.Range("J10").Value = 0.4144845655069451
$ws.Range("O10").Value = 0.876943120004562
$ws.Range("P10").Value = 0.876943120004562
$ws.Range("S10").Value = 0.3634793880693958
$ws.Range("T10").Value = 0.3634793880693957
